$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.840.26'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").Value = '1.864.59'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''247.29'
$ws.Range("E5").Value = '  +2.32%  '
$ws.Range("D6").Value = '''0.6400'
$ws.Range("E6").Value = '  +3.79%  '
$ws.Range("D7").Value = '''1.0000'
$ws.Range("D8").Value = '''0.3016'
$ws.Range("D9").Value = '''0.07515'
$ws.Range("E9").Value = '  +2.40%  '
$ws.Range("D10").Value = '''24.29'
$ws.Range("E10").Value = '  +5.88%  '
$ws.Range("D11").Value = '''0.07679'
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '1.839.59'
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '''5.078'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").Value = '''0.6926'
$ws.Range("E14").Value = '  +4.69%  '
$ws.Range("D15").Value = '''84.77'
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("D16").Value = '''0.000009530'
$ws.Range("E16").Value = '  +6.49%  '
$ws.Range("D17").Value = '''6.135'
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").Value = '29.823.84'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = '2.098.15'
$ws.Range("E19").Value = '  +2.33%  '
$ws.Range("D20").Value = '''240.87'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").Value = '''0.9996'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '''7.415'
$ws.Range("E23").Value = '  +3.93%  '
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '''159.67'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").Value = '''0.1434'
$ws.Range("E26").Value = '  +1.68%  '
$ws.Range("D27").Value = '''8.576'
$ws.Range("E27").Value = '  +1.65%  '
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("D29").Value = '''1.512'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = '''0.06032'
$ws.Range("E30").Value = '  +7.93%  '
$ws.Range("D31").Value = '''1.267'
$ws.Range("E31").Value = '  +5.18%  '
$ws.Range("D32").Value = '''4.152'
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("D33").Value = '''4.156'
$ws.Range("E33").Value = '  +1.30%  '
$ws.Range("D34").Value = '''1.875'
$ws.Range("E34").Value = '  +2.69%  '
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("D36").Value = '''0.7372'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").Value = '''2.621'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = '''2.878'
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("D39").Value = '1.230.21'
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").Value = '''0.01795'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = '''6.412'
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").Value = '''0.9249'
$ws.Range("E42").Value = '  +3.45%  '
$ws.Range("D43").Value = '''1.001'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '2.023.27'
$ws.Range("E44").Value = '  +3.34%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '''102.59'
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("D46").Value = '''66.74'
$ws.Range("E46").Value = '  +3.28%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.00000000122'
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.5086'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''9.367'
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("D50").Value = '''0.4110'
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").Value = '''0.1148'
$ws.Range("E51").Value = '  +3.19%  '
